{"js": "// Office.js (Word JavaScript API) version of the edit.\n//\n// Target change (per the supplied diff):\n//   1. Delete the 2nd paragraph entirely\n//      (\"This is the document you will need to change.  Delete\n//      everything below the above line.\").\n//   2. Replace the text of the 3rd paragraph (formerly \"Add\n//      instructions for your tutor to pull (merge) your request to\n//      the mainline.  Note that the tutor will not merge all\n//      requests, but you should say how it could be achieved.\")\n//      with new instructions about merging a pull request, while\n//      keeping the existing \"_GoBack\" bookmark at the end of the\n//      paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Step 1: drop the \"This is the document...\" paragraph (index 1) ---\nparagraphs.items[1].delete();\nawait context.sync();\n\n// --- Step 2: rewrite the final paragraph's text ---\n// (Re-load paragraphs since the collection shifted after the delete.)\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst target = paragraphs.items[1];\n\nconst newText =\n  \"Within the \\u201cPull requests\\u201d tab in the folder repository, \" +\n  \"select the pull request you would like to merge. \" +\n  \"Click \\u201cMerge pull request\\u201d, and then \\u201cConfirm merge\\u201d \" +\n  \"to commit the document merge (you may add comments if you wish). \" +\n  \"If the document has no merge conflicts with the main branch, it will be merged.\";\n\n// Replace the whole paragraph's content (this also removes the old\n// \"_GoBack\" bookmark, since it sat inside the replaced range).\ntarget.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-create the \"_GoBack\" bookmark, collapsed at the end of the\n// paragraph, matching its original (empty) position in the source.\nconst endRange = target.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) version of the edit.\n#\n# Target change (per the supplied diff):\n#   1. Delete the 2nd paragraph entirely\n#      (\"This is the document you will need to change.  Delete\n#      everything below the above line.\").\n#   2. Replace the text of the 3rd paragraph (formerly \"Add\n#      instructions for your tutor to pull (merge) your request to\n#      the mainline.  Note that the tutor will not merge all\n#      requests, but you should say how it could be achieved.\")\n#      with new instructions about merging a pull request, while\n#      keeping the existing \"_GoBack\" bookmark at the end of the\n#      paragraph.\n\n$d = $word.ActiveDocument\n\n$newText = \"Within the \" + [char]0x201C + \"Pull requests\" + [char]0x201D + `\n    \" tab in the folder repository, select the pull request you would like to merge. \" + `\n    \"Click \" + [char]0x201C + \"Merge pull request\" + [char]0x201D + \", and then \" + `\n    [char]0x201C + \"Confirm merge\" + [char]0x201D + \" to commit the document merge \" + `\n    \"(you may add comments if you wish). If the document has no merge conflicts \" + `\n    \"with the main branch, it will be merged.\"\n\n# The 3rd paragraph already carries a collapsed \"_GoBack\" bookmark in\n# the middle of its text. Re-use that same bookmark object (instead of\n# deleting + recreating it) so it keeps tracking a position inside the\n# document rather than ever landing exactly on the end-of-story mark.\n\n$bm = $d.Bookmarks(\"_GoBack\")\n\n# Delete everything in paragraph 3 that comes AFTER the bookmark.\n$p3Range = $d.Paragraphs(3).Range\n[void]$p3Range.MoveEnd(1, -1)   # exclude the trailing paragraph mark\n$afterBm = $d.Range($bm.Range.End, $p3Range.End)\nif ($afterBm.Start -lt $afterBm.End) {\n    $afterBm.Delete()\n}\n\n# Delete everything in paragraph 3 that comes BEFORE the bookmark.\n$bm = $d.Bookmarks(\"_GoBack\")\n$p3Start = $d.Paragraphs(3).Range.Start\n$beforeBm = $d.Range($p3Start, $bm.Range.Start)\nif ($beforeBm.Start -lt $beforeBm.End) {\n    $beforeBm.Delete()\n}\n\n# Insert the new text immediately before the (now collapsed, alone in\n# its paragraph) bookmark -- this leaves the bookmark collapsed right\n# after the inserted text, exactly like the source XML.\n$bm = $d.Bookmarks(\"_GoBack\")\n$insertPoint = $d.Range($bm.Range.Start, $bm.Range.Start)\n$insertPoint.InsertBefore($newText)\n\n# Finally, delete the \"This is the document you will need to\n# change...\" paragraph (paragraph 2) entirely.\n$d.Paragraphs(2).Range.Delete()\n"}
